$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fifa_world_cup_2018_matches")

# Update match statuses to "completed" and fill in final scores
# Row 38 (match 37): score stays 0-0, just status changes
$ws.Range("D38").Value = "completed"

# Row 39 (match 38): away goals 0 -> 2
$ws.Range("D39").Value = "completed"
$ws.Range("H39").Value = 2

# Row 40 (match 39): home goals 0 -> 1, away goals 0 -> 2
$ws.Range("D40").Value = "completed"
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 2

# Row 41 (match 40): home goals 0 -> 1, away goals 0 -> 2
$ws.Range("D41").Value = "completed"
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 2

# Row 42 (match 41): away goals 0 -> 3
$ws.Range("D42").Value = "completed"
$ws.Range("H42").Value = 3

# Row 43 (match 42): home goals 0 -> 2
$ws.Range("D43").Value = "completed"
$ws.Range("G43").Value = 2

# Row 44 (match 43): away goals 0 -> 2
$ws.Range("D44").Value = "completed"
$ws.Range("H44").Value = 2

# Row 45 (match 44): home goals 0 -> 2, away goals 0 -> 2
$ws.Range("D45").Value = "completed"
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 2

# Update the active cell selection to reflect the final edited cell
$ws.Range("F46").Select()

$wb.Save()
